$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "RXNO_DEF" header in F1, reusing the same formatting (bold, centered,
# thin border) as the other header cells by copying E1's formatting over.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "RXNO_DEF"

# Fill F2:F13 with the literal string "[]" for each data row
for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = "[]"
}
